# Add dynamically-created form fields as a one-row header export: one
# column per field, written as a simple header record in row 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Asset Identification Grp-Field2",
    "Asset Identification Grp-Field1",
    "Asset Identification Sub-Field1",
    "Asset Identification Sub-Field2",
    "Asset Identification Sub-Field3",
    "Asset Identification Sub-Field4"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Widen each field's column so the header text is fully visible.
# NOTE: the COM layer stores ColumnWidth with a constant +5/6 character
# padding baked in (e.g. setting 40 round-trips as 40.8333...), so we
# dial the input back by 5/6 to land on an on-disk width of exactly 40.
$ws.Range("A1:F1").ColumnWidth = 40 - (5/6)
